$wb = $excel.ActiveWorkbook

# This script applies updated market-price / profit figures scraped by the
# scheduled runner into the per-job "Chocobo_Profits" sheets (ALC, ARM, BSM,
# CRP, CUL, GSM, LTW, WVR). All cells are plain literal numbers (no formulas).

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H75").Value = 38457.5
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 38457.5
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 38457.5
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -40329.5

$ws.Range("H78").Value = 38457.5
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 38457.5
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 115372.5
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -124732.5

$ws.Range("H123").Value = 41807.5
$ws.Range("J123").Value = 41807.5
$ws.Range("L123").Value = 41807.5
$ws.Range("N123").Value = -51607.5

$ws.Range("H137").Value = 4942.5
$ws.Range("I137").Value = 2555.7144
$ws.Range("K137").Value = 7667.1432
$ws.Range("M137").Value = -5117.1432

$ws.Range("H138").Value = 1759.6768
$ws.Range("I138").Value = 606.35895
$ws.Range("J138").Value = 2509.3333
$ws.Range("K138").Value = 1819.07685
$ws.Range("L138").Value = 7527.999899999999
$ws.Range("M138").Value = 3320.92315
$ws.Range("N138").Value = -17807.9999


# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 747
$ws.Range("I2").Value = 634.6957
$ws.Range("K2").Value = 634.6957
$ws.Range("M2").Value = -521.6957

$ws.Range("H32").Value = 5459.7744
$ws.Range("I32").Value = 4862.875
$ws.Range("J32").Value = 7506.2856
$ws.Range("K32").Value = 4862.875
$ws.Range("L32").Value = 7506.2856
$ws.Range("M32").Value = -4575.875
$ws.Range("N32").Value = -8080.2856

$ws.Range("H97").Value = 793.13043
$ws.Range("I97").Value = 752.35297
$ws.Range("J97").Value = 908.6667
$ws.Range("K97").Value = 752.35297
$ws.Range("L97").Value = 908.6667
$ws.Range("M97").Value = -256.35297
$ws.Range("N97").Value = -1900.6667

$ws.Range("H116").Value = 747
$ws.Range("I116").Value = 634.6957
$ws.Range("K116").Value = 634.6957
$ws.Range("M116").Value = 1659.3043


# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 747
$ws.Range("I3").Value = 634.6957
$ws.Range("K3").Value = 634.6957
$ws.Range("M3").Value = -520.6957

$ws.Range("H53").Value = 26853.334
$ws.Range("J53").Value = 26853.334
$ws.Range("L53").Value = 26853.334
$ws.Range("N53").Value = -28001.334


# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 6537077.5
$ws.Range("I16").Value = 9260211
$ws.Range("J16").Value = 1556.6
$ws.Range("K16").Value = 9260211
$ws.Range("L16").Value = 1556.6
$ws.Range("M16").Value = -9259924
$ws.Range("N16").Value = -2130.6

$ws.Range("H31").Value = 17244092
$ws.Range("I31").Value = 1023.7857
$ws.Range("J31").Value = 33337622
$ws.Range("K31").Value = 1023.7857
$ws.Range("L31").Value = 33337622
$ws.Range("M31").Value = -728.7857
$ws.Range("N31").Value = -33338212

$ws.Range("H34").Value = 17244092
$ws.Range("I34").Value = 1023.7857
$ws.Range("J34").Value = 33337622
$ws.Range("K34").Value = 1023.7857
$ws.Range("L34").Value = 33337622
$ws.Range("M34").Value = -821.7857
$ws.Range("N34").Value = -33338026

$ws.Range("H106").Value = 32800
$ws.Range("J106").Value = 32800
$ws.Range("L106").Value = 32800
$ws.Range("N106").Value = -35324

$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws.Range("H113").Value = 6537077.5
$ws.Range("I113").Value = 9260211
$ws.Range("J113").Value = 1556.6
$ws.Range("K113").Value = 9260211
$ws.Range("L113").Value = 1556.6
$ws.Range("M113").Value = -9258041
$ws.Range("N113").Value = -5896.6

$ws.Range("H132").Value = 2620.327
$ws.Range("I132").Value = 2122.0256
$ws.Range("J132").Value = 4115.231
$ws.Range("K132").Value = 6366.0768
$ws.Range("L132").Value = 12345.693
$ws.Range("M132").Value = -3836.0768
$ws.Range("N132").Value = -17405.693


# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 164.57143
$ws.Range("I33").Value = 191
$ws.Range("J33").Value = 98.5
$ws.Range("K33").Value = 1146
$ws.Range("L33").Value = 591
$ws.Range("M33").Value = -863
$ws.Range("N33").Value = -1157

$ws.Range("H94").Value = 2373.5
$ws.Range("I94").Value = 845
$ws.Range("J94").Value = 2883
$ws.Range("K94").Value = 2535
$ws.Range("L94").Value = 8649
$ws.Range("M94").Value = -1859
$ws.Range("N94").Value = -10001

$ws.Range("H96").Value = 6967
$ws.Range("J96").Value = 6967
$ws.Range("L96").Value = 20901
$ws.Range("N96").Value = -25019

$ws.Range("H113").Value = 734.0303
$ws.Range("I113").Value = 606.55554
$ws.Range("K113").Value = 1819.66662
$ws.Range("M113").Value = 350.33338

$ws.Range("H131").Value = 6667494.5
$ws.Range("I131").Value = 100000280
$ws.Range("J131").Value = 867.0571
$ws.Range("K131").Value = 300000840
$ws.Range("L131").Value = 2601.1713
$ws.Range("M131").Value = -299995800
$ws.Range("N131").Value = -12681.1713

$ws.Range("H133").Value = 2939.375
$ws.Range("J133").Value = 2960
$ws.Range("L133").Value = 8880
$ws.Range("N133").Value = -19000


# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 8860169
$ws.Range("I14").Value = 10333364
$ws.Range("J14").Value = 20999
$ws.Range("K14").Value = 10333364
$ws.Range("L14").Value = 20999
$ws.Range("M14").Value = -10333196
$ws.Range("N14").Value = -21335


# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4091.851
$ws.Range("I132").Value = 1760.2812
$ws.Range("J132").Value = 9065.866
$ws.Range("K132").Value = 5280.8436
$ws.Range("L132").Value = 27197.598
$ws.Range("M132").Value = -2750.8436
$ws.Range("N132").Value = -32257.598

$ws.Range("H134").Value = 39802.375
$ws.Range("J134").Value = 39802.375
$ws.Range("L134").Value = 39802.375
$ws.Range("N134").Value = -49942.375


# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 22965
$ws.Range("I41").Value = 3342
$ws.Range("J41").Value = 32776.5
$ws.Range("K41").Value = 3342
$ws.Range("L41").Value = 32776.5
$ws.Range("M41").Value = -2952
$ws.Range("N41").Value = -33556.5

$ws.Range("H45").Value = 5833.3335
$ws.Range("J45").Value = 6250
$ws.Range("L45").Value = 6250
$ws.Range("N45").Value = -7232

$ws.Range("H74").Value = 5417.125
$ws.Range("I74").Value = 5486.3335
$ws.Range("J74").Value = 5375.6
$ws.Range("K74").Value = 5486.3335
$ws.Range("L74").Value = 5375.6
$ws.Range("M74").Value = -4550.3335
$ws.Range("N74").Value = -7247.6

$ws.Range("H77").Value = 5417.125
$ws.Range("I77").Value = 5486.3335
$ws.Range("J77").Value = 5375.6
$ws.Range("K77").Value = 16459.0005
$ws.Range("L77").Value = 16126.8
$ws.Range("M77").Value = -11779.0005
$ws.Range("N77").Value = -25486.8

$ws.Range("H100").Value = 421.375
$ws.Range("I100").Value = 424.42856
$ws.Range("J100").Value = 400
$ws.Range("K100").Value = 848.85712
$ws.Range("L100").Value = 800
$ws.Range("M100").Value = -307.85712
$ws.Range("N100").Value = -1882

$ws.Range("H136").Value = 2400.8948
$ws.Range("I136").Value = 678.13336
$ws.Range("J136").Value = 8861.25
$ws.Range("K136").Value = 2034.40008
$ws.Range("L136").Value = 26583.75
$ws.Range("M136").Value = 515.5999199999999
$ws.Range("N136").Value = -31683.75

$ws.Range("H140").Value = 38818.918
$ws.Range("J140").Value = 38818.918
$ws.Range("L140").Value = 38818.918
$ws.Range("N140").Value = -49178.918


Write-Output "Updated market data across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets."